# "Añadiendo nombre y apellidos"
#
# 1. Remove the stray _GoBack bookmark that currently sits at the end of
#    the "...hora actual (comitea)." list item.
# 2. Append a new (non-numbered) paragraph at the very end of the body,
#    styled like the other list items (Prrafodelista / left indent 426
#    twips, but without the bullet numbering), containing the author's
#    name, a lastRenderedPageBreak marker before the text (as produced by
#    Word when it repaginates), and a fresh _GoBack bookmark collapsed at
#    the end of the new text.

$d = $word.ActiveDocument

# --- Step 1: drop the old _GoBack bookmark -------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: append the new paragraph -------------------------------------
# Collapse a Range sitting at the very end of the document's main story and
# inject a fully-formed WordprocessingML paragraph through InsertXML so we
# get full control over pPr / run contents (incl. lastRenderedPageBreak and
# the bookmark) in one shot.
$endRange = $d.Content
$endRange.Collapse(0)

$eAcute = [char]0x00E9
$aAcute = [char]0x00E1
$name = "N" + $eAcute + "stor Rold" + $aAcute + "n Aznar"

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="426"/></w:pPr>' +
       '<w:r><w:lastRenderedPageBreak/><w:t>' + $name + '</w:t></w:r>' +
       '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
       '</w:p>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$endRange.InsertXML($xml)
